$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace formula-driven values in B2:B5 with literal numeric values (formulas removed)
$ws.Range("B2").Value = [double]"1.30946032491679E-2"
$ws.Range("B3").Value = [double]"1.1538744957234099E-5"
$ws.Range("B4").Value = [double]"7.5039162349632205E-8"
$ws.Range("B5").Value = [double]"2.7127485839104799E-9"

# Clear the cycNo values (A6:A11) and the remaining shared-formula results (B6:B11)
$ws.Range("A6:A11").ClearContents()
$ws.Range("B6:B11").ClearContents()

# Remove the now-unused trailing rows 25:30, shrinking the sheet's used range
[void]$ws.Range("A25:B30").Delete(-4162)

# Move/update the active selection shown in the sheet view
[void]$ws.Range("E8").Select()
